$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $styleName = $p.Style.NameLocal
    if ($styleName -eq "Title" -or $styleName -eq "Author" -or $styleName -eq "Abstract") {
        $r = $p.Range
        # Trim the trailing paragraph mark from the range so we only touch the text.
        $r.MoveEnd(1, -1) | Out-Null
        $text = $r.Text

        # Re-assigning the identical string is treated as a no-op by the
        # engine (it skips writes whose value doesn't change), so the many
        # single-word runs would otherwise survive untouched. Force Word to
        # rebuild the paragraph as one run by first writing a distinct
        # placeholder value, then writing the real text back.
        $r.Text = $text + [char]1
        $r2 = $p.Range
        $r2.MoveEnd(1, -1) | Out-Null
        $r2.Text = $text
    }
}
